# Add the new "web" command dragAndDrop(fromLocator,toLocator) to the
# alphabetically-sorted list kept in column U of the hidden '#system' sheet.
#
# The list lives in U2:U111 (named range "web"). The new entry belongs
# between "doubleClickByLabelAndWait(label,waitMs)" (U58) and
# "editLocalStorage(key,value)" (U59), so every entry from U59 downward
# needs to shift down by one row, the new command goes into U59, and the
# named range + sheet dimension need to grow to U112 / row 112.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$insertRow = 59
$lastRow = 111
$col = 21  # column U

# Shift existing values down by one row, starting from the bottom so we
# never overwrite a value before it has been copied.
for ($r = $lastRow; $r -ge $insertRow; $r--) {
    $source = $ws.Cells.Item($r, $col)
    $target = $ws.Cells.Item($r + 1, $col)
    $target.Value = $source.Value2
}

# Insert the new command text at its alphabetically-correct position.
$ws.Cells.Item($insertRow, $col).Value = "dragAndDrop(fromLocator,toLocator)"

# Grow the "web" defined name to cover the extra row.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "web") {
        $n.RefersTo = "='#system'!`$U`$2:`$U`$112"
    }
}
